# Update document date and multiplication problems per diff.
$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-25 Saturday", "2023-11-26 Sunday"),
    @("80×26=", "48×38="),
    @("75×43=", "54×19="),
    @("87×55=", "45×99="),
    @("56×14=", "12×57="),
    @("33×60=", "54×73="),
    @("75×14=", "62×80="),
    @("88×23=", "67×31="),
    @("75×28=", "84×47="),
    @("24×62=", "71×24="),
    @("13×14=", "21×90="),
    @("18×58=", "43×72="),
    @("54×54=", "26×80="),
    @("67×38=", "55×62="),
    @("87×13=", "70×40="),
    @("39×69=", "50×72="),
    @("66×48=", "34×71="),
    @("83×52=", "16×71="),
    @("88×55=", "67×95="),
    @("38×57=", "50×71="),
    @("60×57=", "16×92="),
    @("41×50=", "42×52="),
    @("99×43=", "33×27="),
    @("51×42=", "19×26="),
    @("47×71=", "21×88="),
    @("18×95=", "41×36=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done applying replacements."
